# Apply the "log and config files" edits to Sheet1:
#  - new numeric log values in column G/J (one of them a formula)
#  - a new formula in P6
#  - a new percentage formula in A9 formatted with the built-in "Percent" cell style
#  - move the active selection to D4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G1: precision-ish ratio, entered as a formula
$ws.Range("G1").Formula = "=0.5306/0.984"

# J1 / J2 / J5: plain numeric log values
$ws.Range("J1").Value = 200.21
$ws.Range("J2").Value = -142.17
$ws.Range("J5").Value = 58.05

# P6: simple arithmetic formula
$ws.Range("P6").Formula = "=12*5+10"

# A9: recomputed recall-style formula, formatted as a percentage via the
# built-in "Percent" cell style (adds cellStyleXfs/cellStyles entries)
$ws.Range("A9").Formula = "=1-29762/(29762+34988)"
$ws.Range("A9").Style = "Percent"

# Leave the selection on D4, matching the saved view state
$ws.Range("D4").Select()
